# paises.xlsx refresh: "Update countries & provincias Spain"
# - Refresh the "last updated" timestamp
# - Refresh per-country COVID counters (Casos totales/Nuevos casos/Casos activos/
#   Recuperados/Muertes hoy/Muertes) for the countries whose figures moved
# - A handful of rows change which country they display because Nicaragua, Cuba
#   and Birmania overtook their neighbours in the ranking (sorted by Casos totales)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: refresh the "Datos actualizados ..." timestamp
$ws.Range("A1").Value = "Datos actualizados a 9 de Septiembre de 2020 a las 00:21"

# Row 4: Estados Unidos
$ws.Cells.Item(4, 2).Value = 6507728  # Casos totales
$ws.Cells.Item(4, 3).Value = 22058  # Nuevos casos
$ws.Cells.Item(4, 4).Value = 3785226  # Casos activos
$ws.Cells.Item(4, 5).Value = 2528598  # Recuperados
$ws.Cells.Item(4, 7).Value = 370  # Muertes hoy
$ws.Cells.Item(4, 8).Value = 193904  # Muertes

# Row 6: Brasil
$ws.Cells.Item(6, 2).Value = 4162073  # Casos totales
$ws.Cells.Item(6, 3).Value = 14279  # Nuevos casos
$ws.Cells.Item(6, 4).Value = 3397234  # Casos activos
$ws.Cells.Item(6, 5).Value = 637375  # Recuperados
$ws.Cells.Item(6, 7).Value = 463  # Muertes hoy
$ws.Cells.Item(6, 8).Value = 127464  # Muertes

# Row 9: Colombia
$ws.Cells.Item(9, 2).Value = 679513  # Casos totales
$ws.Cells.Item(9, 3).Value = 7665  # Nuevos casos
$ws.Cells.Item(9, 4).Value = 541462  # Casos activos
$ws.Cells.Item(9, 5).Value = 116234  # Recuperados
$ws.Cells.Item(9, 7).Value = 202  # Muertes hoy
$ws.Cells.Item(9, 8).Value = 21817  # Muertes

# Row 29: Canada
$ws.Cells.Item(29, 2).Value = 133748  # Casos totales
$ws.Cells.Item(29, 3).Value = 1606  # Nuevos casos
$ws.Cells.Item(29, 4).Value = 117565  # Casos activos
$ws.Cells.Item(29, 5).Value = 7030  # Recuperados
$ws.Cells.Item(29, 7).Value = 7  # Muertes hoy
$ws.Cells.Item(29, 8).Value = 9153  # Muertes

# Row 43: Guatemala
$ws.Cells.Item(43, 2).Value = 78721  # Casos totales
$ws.Cells.Item(43, 3).Value = 893  # Nuevos casos
$ws.Cells.Item(43, 4).Value = 67462  # Casos activos
$ws.Cells.Item(43, 5).Value = 8369  # Recuperados
$ws.Cells.Item(43, 7).Value = 28  # Muertes hoy
$ws.Cells.Item(43, 8).Value = 2890  # Muertes

# Row 48: Japon
$ws.Cells.Item(48, 2).Value = 72234  # Casos totales
$ws.Cells.Item(48, 3).Value = 378  # Nuevos casos
$ws.Cells.Item(48, 4).Value = 63282  # Casos activos
$ws.Cells.Item(48, 5).Value = 7575  # Recuperados
$ws.Cells.Item(48, 7).Value = 14  # Muertes hoy
$ws.Cells.Item(48, 8).Value = 1377  # Muertes

# Row 54: Barein
$ws.Cells.Item(54, 2).Value = 56778  # Casos totales
$ws.Cells.Item(54, 3).Value = 702  # Nuevos casos
$ws.Cells.Item(54, 4).Value = 51574  # Casos activos
$ws.Cells.Item(54, 5).Value = 5002  # Recuperados

# Row 84: Bulgaria
$ws.Cells.Item(84, 2).Value = 17313  # Casos totales
$ws.Cells.Item(84, 3).Value = 167  # Nuevos casos
$ws.Cells.Item(84, 4).Value = 12297  # Casos activos
$ws.Cells.Item(84, 5).Value = 4324  # Recuperados
$ws.Cells.Item(84, 7).Value = 15  # Muertes hoy
$ws.Cells.Item(84, 8).Value = 692  # Muertes

# Row 89: Zambia
$ws.Cells.Item(89, 2).Value = 12952  # Casos totales
$ws.Cells.Item(89, 3).Value = 116  # Nuevos casos
$ws.Cells.Item(89, 4).Value = 11787  # Casos activos
$ws.Cells.Item(89, 5).Value = 868  # Recuperados
$ws.Cells.Item(89, 7).Value = 2  # Muertes hoy
$ws.Cells.Item(89, 8).Value = 297  # Muertes

# Row 97: Guayana Francesa
$ws.Cells.Item(97, 2).Value = 9387  # Casos totales
$ws.Cells.Item(97, 3).Value = 32  # Nuevos casos
$ws.Cells.Item(97, 4).Value = 8946  # Casos activos
$ws.Cells.Item(97, 5).Value = 379  # Recuperados

# Row 105: Zimbabue
$ws.Cells.Item(105, 2).Value = 7388  # Casos totales
$ws.Cells.Item(105, 3).Value = 90  # Nuevos casos
$ws.Cells.Item(105, 4).Value = 5477  # Casos activos
$ws.Cells.Item(105, 5).Value = 1693  # Recuperados
$ws.Cells.Item(105, 7).Value = 8  # Muertes hoy
$ws.Cells.Item(105, 8).Value = 218  # Muertes

# Row 107: Luxemburgo
$ws.Cells.Item(107, 2).Value = 6974  # Casos totales
$ws.Cells.Item(107, 3).Value = 14  # Nuevos casos
$ws.Cells.Item(107, 5).Value = 594  # Recuperados

# Row 109: Malaui
$ws.Cells.Item(109, 2).Value = 5630  # Casos totales
$ws.Cells.Item(109, 3).Value = 9  # Nuevos casos
$ws.Cells.Item(109, 4).Value = 3630  # Casos activos
$ws.Cells.Item(109, 5).Value = 1824  # Recuperados

# Row 116: Nicaragua
$ws.Cells.Item(116, 1).Value = "Nicaragua"  # Pais
$ws.Cells.Item(116, 2).Value = 4818  # Casos totales
$ws.Cells.Item(116, 3).Value = 150  # Nuevos casos
$ws.Cells.Item(116, 4).Value = 2913  # Casos activos
$ws.Cells.Item(116, 5).Value = 1761  # Recuperados
$ws.Cells.Item(116, 7).Value = 3  # Muertes hoy
$ws.Cells.Item(116, 8).Value = 144  # Muertes

# Row 117: Republica de Africa Central
$ws.Cells.Item(117, 1).Value = "Republica de Africa Central"  # Pais
$ws.Cells.Item(117, 2).Value = 4735  # Casos totales
$ws.Cells.Item(117, 3).Value = 6  # Nuevos casos
$ws.Cells.Item(117, 4).Value = 1825  # Casos activos
$ws.Cells.Item(117, 5).Value = 2848  # Recuperados
$ws.Cells.Item(117, 8).Value = 62  # Muertes

# Row 118: Eslovaquia
$ws.Cells.Item(118, 1).Value = "Eslovaquia"  # Pais
$ws.Cells.Item(118, 2).Value = 4727  # Casos totales
$ws.Cells.Item(118, 3).Value = 91  # Nuevos casos
$ws.Cells.Item(118, 5).Value = 1777  # Recuperados
$ws.Cells.Item(118, 8).Value = 37  # Muertes

# Row 120: Ruanda
$ws.Cells.Item(120, 2).Value = 4439  # Casos totales
$ws.Cells.Item(120, 3).Value = 30  # Nuevos casos
$ws.Cells.Item(120, 4).Value = 2307  # Casos activos
$ws.Cells.Item(120, 5).Value = 2112  # Recuperados
$ws.Cells.Item(120, 7).Value = 1  # Muertes hoy
$ws.Cells.Item(120, 8).Value = 20  # Muertes

# Row 122: Cuba
$ws.Cells.Item(122, 1).Value = "Cuba"  # Pais
$ws.Cells.Item(122, 2).Value = 4377  # Casos totales
$ws.Cells.Item(122, 3).Value = 25  # Nuevos casos
$ws.Cells.Item(122, 4).Value = 3700  # Casos activos
$ws.Cells.Item(122, 5).Value = 573  # Recuperados
$ws.Cells.Item(122, 7).Value = 2  # Muertes hoy
$ws.Cells.Item(122, 8).Value = 104  # Muertes

# Row 123: Surinam
$ws.Cells.Item(123, 1).Value = "Surinam"  # Pais
$ws.Cells.Item(123, 2).Value = 4360  # Casos totales
$ws.Cells.Item(123, 4).Value = 3544  # Casos activos
$ws.Cells.Item(123, 5).Value = 725  # Recuperados
$ws.Cells.Item(123, 8).Value = 91  # Muertes

# Row 152: Birmania
$ws.Cells.Item(152, 1).Value = "Birmania"  # Pais
$ws.Cells.Item(152, 2).Value = 1709  # Casos totales
$ws.Cells.Item(152, 3).Value = 191  # Nuevos casos
$ws.Cells.Item(152, 4).Value = 460  # Casos activos
$ws.Cells.Item(152, 5).Value = 1239  # Recuperados
$ws.Cells.Item(152, 7).Value = 2  # Muertes hoy
$ws.Cells.Item(152, 8).Value = 10  # Muertes

# Row 153: Uruguay
$ws.Cells.Item(153, 1).Value = "Uruguay"  # Pais
$ws.Cells.Item(153, 2).Value = 1693  # Casos totales
$ws.Cells.Item(153, 3).Value = 0  # Nuevos casos
$ws.Cells.Item(153, 4).Value = 1466  # Casos activos
$ws.Cells.Item(153, 5).Value = 182  # Recuperados
$ws.Cells.Item(153, 8).Value = 45  # Muertes

# Row 156: Republica de Chipre
$ws.Cells.Item(156, 2).Value = 1511  # Casos totales
$ws.Cells.Item(156, 3).Value = 1  # Nuevos casos
$ws.Cells.Item(156, 5).Value = 252  # Recuperados

# Row 178: Comoras
$ws.Cells.Item(178, 2).Value = 456  # Casos totales
$ws.Cells.Item(178, 3).Value = 4  # Nuevos casos
$ws.Cells.Item(178, 4).Value = 415  # Casos activos
$ws.Cells.Item(178, 5).Value = 34  # Recuperados
